$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 213, pushing the existing rows 213-304 down to 214-305
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new data record
$ws.Cells.Item(213, 1).Value = 3
$ws.Cells.Item(213, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = 44726
$ws.Cells.Item(213, 5).Value = 5
$ws.Cells.Item(213, 6).Value = 100112001
$ws.Cells.Item(213, 7).Value = "Berenjena"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 105
$ws.Cells.Item(213, 11).Value = 5500
$ws.Cells.Item(213, 12).Value = 6000
$ws.Cells.Item(213, 13).Value = 5738
$ws.Cells.Item(213, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(213, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(213, 16).Value = 96
$ws.Cells.Item(213, 17).Value = 60
$ws.Cells.Item(213, 18).Value = "Hortaliza"
